$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.172.26"
$ws.Range("E2").Value = "  -1.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.91"
$ws.Range("E3").Value = "  -3.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.84"

# Row 6
$ws.Range("E6").Value = "  +0.19%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("E7").Value = "  -2.34%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2819"
$ws.Range("E8").Value = "  -1.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06550"
$ws.Range("E9").Value = "  -3.21%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.08"
$ws.Range("E10").Value = "  +2.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07828"
$ws.Range("E11").Value = "  +0.62%  "

# Row 12
$ws.Range("E12").Value = "  -7.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.99"
$ws.Range("E13").Value = "  -2.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.109"
$ws.Range("E14").Value = "  -2.77%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6658"
$ws.Range("E15").Value = "  -1.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.87"
$ws.Range("E16").Value = "  -3.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.217.45"
$ws.Range("E17").Value = "  -1.23%  "

# Row 18
$ws.Range("E18").Value = "  +0.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.453"
$ws.Range("E19").Value = "  -0.46%  "

# Row 20
$ws.Range("E20").Value = "  -1.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.107.76"
$ws.Range("E21").Value = "  -2.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007242"
$ws.Range("E22").Value = "  -3.98%  "

# Row 23
$ws.Range("E23").Value = "  +0.14%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.139"
$ws.Range("E24").Value = "  -3.58%  "

# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.320"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.53"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  -3.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920"
$ws.Range("E28").Value = "  -9.18%  "

# Row 29
$ws.Range("E29").Value = "  -3.46%  "

# Row 30
$ws.Range("E30").Value = "  -4.37%  "

# Row 31
$ws.Range("E31").Value = "  -3.95%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.470"
$ws.Range("E32").Value = "  -3.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.094"
$ws.Range("E33").Value = "  -4.68%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04654"
$ws.Range("E34").Value = "  -2.44%  "

# Row 35
$ws.Range("E35").Value = "  -1.76%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6995"
$ws.Range("E36").Value = "  -4.52%  "

# Row 37
$ws.Range("E37").Value = "  +0.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("E38").Value = "  -0.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -3.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.386"
$ws.Range("E40").Value = "  -0.16%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.513"
$ws.Range("E41").Value = "  -4.26%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.02"
$ws.Range("E42").Value = "  -4.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8534"
$ws.Range("E43").Value = "  -1.23%  "

# Row 44
$ws.Range("E44").Value = "  -2.81%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4163"
$ws.Range("E46").Value = "  -2.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.68"
$ws.Range("E47").Value = "  -2.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "997.58"
$ws.Range("E48").Value = "  +2.18%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.194"
$ws.Range("E49").Value = "  -3.69%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.087"
$ws.Range("E50").Value = "  +2.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.92"
$ws.Range("E51").Value = "  -2.84%  "
